# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.950.17"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.649.11"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.59"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.82"
$ws.Range("E10").Value = "  +4.53%  "
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.880.55"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "1.654.35"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.70"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").Value = "26.969.46"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.78"
$ws.Range("E19").Value = "  +4.60%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.40"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.65"
$ws.Range("E22").Value = "  +7.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.43"
$ws.Range("E23").Value = "  +4.46%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.26"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.39"
$ws.Range("E27").Value = "  +4.95%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.91"
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0511"
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").Value = "1.250.39"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("E39").Value = "  +3.26%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").Value = "1.792.78"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  -4.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.35"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.57"
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.60"
$ws.Range("E51").Value = "  +0.85%  "
